# Apply NATMI TPM re-run update to the Wnt5a-Fzd7 LR-pair sheet.
# Rows 2-7: existing "FAPs -> *" rows get refreshed TPM-derived values
#           and the sending cluster is corrected to "ECs".
# Rows 8-13: newly added "FAPs -> *" rows (original FAPs data, shifted down).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: ECs -> ECs
$ws.Cells.Item(2, 1).Value = "ECs"
$ws.Cells.Item(2, 2).Value = "Wnt5a"
$ws.Cells.Item(2, 3).Value = "Fzd7"
$ws.Cells.Item(2, 4).Value = "ECs"
$ws.Cells.Item(2, 5).Value = 2
$ws.Cells.Item(2, 6).Value = 1
$ws.Cells.Item(2, 7).Value = 0.125615
$ws.Cells.Item(2, 8).Value = 0.25123
$ws.Cells.Item(2, 9).Value = 0.02647478672532295
$ws.Cells.Item(2, 10).Value = 0.01780700335556722
$ws.Cells.Item(2, 11).Value = 2
$ws.Cells.Item(2, 12).Value = 1
$ws.Cells.Item(2, 13).Value = 1.905108
$ws.Cells.Item(2, 14).Value = 3.810216
$ws.Cells.Item(2, 15).Value = 0.07580486173280727
$ws.Cells.Item(2, 16).Value = 0.05715529216076502
$ws.Cells.Item(2, 17).Value = 0.23931014142
$ws.Cells.Item(2, 18).Value = 0.95724056568
$ws.Cells.Item(2, 19).Value = 0.002006917547118668
$ws.Cells.Item(2, 20).Value = 0.001017764479295168

# Row 3: ECs -> FAPs
$ws.Cells.Item(3, 1).Value = "ECs"
$ws.Cells.Item(3, 2).Value = "Wnt5a"
$ws.Cells.Item(3, 3).Value = "Fzd7"
$ws.Cells.Item(3, 4).Value = "FAPs"
$ws.Cells.Item(3, 5).Value = 2
$ws.Cells.Item(3, 6).Value = 1
$ws.Cells.Item(3, 7).Value = 0.125615
$ws.Cells.Item(3, 8).Value = 0.25123
$ws.Cells.Item(3, 9).Value = 0.02647478672532295
$ws.Cells.Item(3, 10).Value = 0.01780700335556722
$ws.Cells.Item(3, 11).Value = 3
$ws.Cells.Item(3, 12).Value = 1
$ws.Cells.Item(3, 13).Value = 3.068283333333333
$ws.Cells.Item(3, 14).Value = 9.20485
$ws.Cells.Item(3, 15).Value = 0.1220879833796353
$ws.Cells.Item(3, 16).Value = 0.1380777076800943
$ws.Cells.Item(3, 17).Value = 0.3854224109166667
$ws.Cells.Item(3, 18).Value = 2.3125344655
$ws.Cells.Item(3, 19).Value = 0.003232253321700619
$ws.Cells.Item(3, 20).Value = 0.002458750203988468

# Row 4: ECs -> Inflammatory-Mac
$ws.Cells.Item(4, 1).Value = "ECs"
$ws.Cells.Item(4, 2).Value = "Wnt5a"
$ws.Cells.Item(4, 3).Value = "Fzd7"
$ws.Cells.Item(4, 4).Value = "Inflammatory-Mac"
$ws.Cells.Item(4, 5).Value = 2
$ws.Cells.Item(4, 6).Value = 1
$ws.Cells.Item(4, 7).Value = 0.125615
$ws.Cells.Item(4, 8).Value = 0.25123
$ws.Cells.Item(4, 9).Value = 0.02647478672532295
$ws.Cells.Item(4, 10).Value = 0.01780700335556722
$ws.Cells.Item(4, 11).Value = 3
$ws.Cells.Item(4, 12).Value = 1
$ws.Cells.Item(4, 13).Value = 3.725954
$ws.Cells.Item(4, 14).Value = 11.177862
$ws.Cells.Item(4, 15).Value = 0.148256911310435
$ws.Cells.Item(4, 16).Value = 0.1676739503331867
$ws.Cells.Item(4, 17).Value = 0.46803571171
$ws.Cells.Item(4, 18).Value = 2.80821427026
$ws.Cells.Item(4, 19).Value = 0.003925070107498887
$ws.Cells.Item(4, 20).Value = 0.002985770596224267

# Row 5: ECs -> MuSCs
$ws.Cells.Item(5, 1).Value = "ECs"
$ws.Cells.Item(5, 2).Value = "Wnt5a"
$ws.Cells.Item(5, 3).Value = "Fzd7"
$ws.Cells.Item(5, 4).Value = "MuSCs"
$ws.Cells.Item(5, 5).Value = 2
$ws.Cells.Item(5, 6).Value = 1
$ws.Cells.Item(5, 7).Value = 0.125615
$ws.Cells.Item(5, 8).Value = 0.25123
$ws.Cells.Item(5, 9).Value = 0.02647478672532295
$ws.Cells.Item(5, 10).Value = 0.01780700335556722
$ws.Cells.Item(5, 11).Value = 2
$ws.Cells.Item(5, 12).Value = 1
$ws.Cells.Item(5, 13).Value = 6.825836
$ws.Cells.Item(5, 14).Value = 13.651672
$ws.Cells.Item(5, 15).Value = 0.2716022158275637
$ws.Cells.Item(5, 16).Value = 0.2047824327132465
$ws.Cells.Item(5, 17).Value = 0.85742738914
$ws.Cells.Item(5, 18).Value = 3.42970955656
$ws.Cells.Item(5, 19).Value = 0.007190610738159884
$ws.Cells.Item(5, 20).Value = 0.003646561466486

# Row 6: ECs -> Neutrophils
$ws.Cells.Item(6, 1).Value = "ECs"
$ws.Cells.Item(6, 2).Value = "Wnt5a"
$ws.Cells.Item(6, 3).Value = "Fzd7"
$ws.Cells.Item(6, 4).Value = "Neutrophils"
$ws.Cells.Item(6, 5).Value = 2
$ws.Cells.Item(6, 6).Value = 1
$ws.Cells.Item(6, 7).Value = 0.125615
$ws.Cells.Item(6, 8).Value = 0.25123
$ws.Cells.Item(6, 9).Value = 0.02647478672532295
$ws.Cells.Item(6, 10).Value = 0.01780700335556722
$ws.Cells.Item(6, 11).Value = 3
$ws.Cells.Item(6, 12).Value = 1
$ws.Cells.Item(6, 13).Value = 4.247626666666666
$ws.Cells.Item(6, 14).Value = 12.74288
$ws.Cells.Item(6, 15).Value = 0.1690144349607748
$ws.Cells.Item(6, 16).Value = 0.1911500632430207
$ws.Cells.Item(6, 17).Value = 0.5335656237333333
$ws.Cells.Item(6, 18).Value = 3.2013937424
$ws.Cells.Item(6, 19).Value = 0.004474621119087479
$ws.Cells.Item(6, 20).Value = 0.003403809817585356

# Row 7: ECs -> Resolving-Mac
$ws.Cells.Item(7, 1).Value = "ECs"
$ws.Cells.Item(7, 2).Value = "Wnt5a"
$ws.Cells.Item(7, 3).Value = "Fzd7"
$ws.Cells.Item(7, 4).Value = "Resolving-Mac"
$ws.Cells.Item(7, 5).Value = 2
$ws.Cells.Item(7, 6).Value = 1
$ws.Cells.Item(7, 7).Value = 0.125615
$ws.Cells.Item(7, 8).Value = 0.25123
$ws.Cells.Item(7, 9).Value = 0.02647478672532295
$ws.Cells.Item(7, 10).Value = 0.01780700335556722
$ws.Cells.Item(7, 11).Value = 3
$ws.Cells.Item(7, 12).Value = 1
$ws.Cells.Item(7, 13).Value = 5.358931000000001
$ws.Cells.Item(7, 14).Value = 16.076793
$ws.Cells.Item(7, 15).Value = 0.213233592788784
$ws.Cells.Item(7, 16).Value = 0.2411605538696867
$ws.Cells.Item(7, 17).Value = 0.6731621175650001
$ws.Cells.Item(7, 18).Value = 4.038972705390001
$ws.Cells.Item(7, 19).Value = 0.005645313891757419
$ws.Cells.Item(7, 20).Value = 0.00429434679198796

# Row 8: FAPs -> ECs
$ws.Cells.Item(8, 1).Value = "FAPs"
$ws.Cells.Item(8, 2).Value = "Wnt5a"
$ws.Cells.Item(8, 3).Value = "Fzd7"
$ws.Cells.Item(8, 4).Value = "ECs"
$ws.Cells.Item(8, 5).Value = 3
$ws.Cells.Item(8, 6).Value = 1
$ws.Cells.Item(8, 7).Value = 4.619088000000001
$ws.Cells.Item(8, 8).Value = 13.857264
$ws.Cells.Item(8, 9).Value = 0.9735252132746771
$ws.Cells.Item(8, 10).Value = 0.9821929966444328
$ws.Cells.Item(8, 11).Value = 2
$ws.Cells.Item(8, 12).Value = 1
$ws.Cells.Item(8, 13).Value = 1.905108
$ws.Cells.Item(8, 14).Value = 3.810216
$ws.Cells.Item(8, 15).Value = 0.07580486173280727
$ws.Cells.Item(8, 16).Value = 0.05715529216076502
$ws.Cells.Item(8, 17).Value = 8.799861501504001
$ws.Cells.Item(8, 18).Value = 52.79916900902401
$ws.Cells.Item(8, 19).Value = 0.07379794418568861
$ws.Cells.Item(8, 20).Value = 0.05613752768146986

# Row 9: FAPs -> FAPs
$ws.Cells.Item(9, 1).Value = "FAPs"
$ws.Cells.Item(9, 2).Value = "Wnt5a"
$ws.Cells.Item(9, 3).Value = "Fzd7"
$ws.Cells.Item(9, 4).Value = "FAPs"
$ws.Cells.Item(9, 5).Value = 3
$ws.Cells.Item(9, 6).Value = 1
$ws.Cells.Item(9, 7).Value = 4.619088000000001
$ws.Cells.Item(9, 8).Value = 13.857264
$ws.Cells.Item(9, 9).Value = 0.9735252132746771
$ws.Cells.Item(9, 10).Value = 0.9821929966444328
$ws.Cells.Item(9, 11).Value = 3
$ws.Cells.Item(9, 12).Value = 1
$ws.Cells.Item(9, 13).Value = 3.068283333333333
$ws.Cells.Item(9, 14).Value = 9.20485
$ws.Cells.Item(9, 15).Value = 0.1220879833796353
$ws.Cells.Item(9, 16).Value = 0.1380777076800943
$ws.Cells.Item(9, 17).Value = 14.1726707256
$ws.Cells.Item(9, 18).Value = 127.5540365304
$ws.Cells.Item(9, 19).Value = 0.1188557300579347
$ws.Cells.Item(9, 20).Value = 0.1356189574761058

# Row 10: FAPs -> Inflammatory-Mac
$ws.Cells.Item(10, 1).Value = "FAPs"
$ws.Cells.Item(10, 2).Value = "Wnt5a"
$ws.Cells.Item(10, 3).Value = "Fzd7"
$ws.Cells.Item(10, 4).Value = "Inflammatory-Mac"
$ws.Cells.Item(10, 5).Value = 3
$ws.Cells.Item(10, 6).Value = 1
$ws.Cells.Item(10, 7).Value = 4.619088000000001
$ws.Cells.Item(10, 8).Value = 13.857264
$ws.Cells.Item(10, 9).Value = 0.9735252132746771
$ws.Cells.Item(10, 10).Value = 0.9821929966444328
$ws.Cells.Item(10, 11).Value = 3
$ws.Cells.Item(10, 12).Value = 1
$ws.Cells.Item(10, 13).Value = 3.725954
$ws.Cells.Item(10, 14).Value = 11.177862
$ws.Cells.Item(10, 15).Value = 0.148256911310435
$ws.Cells.Item(10, 16).Value = 0.1676739503331867
$ws.Cells.Item(10, 17).Value = 17.210509409952
$ws.Cells.Item(10, 18).Value = 154.894584689568
$ws.Cells.Item(10, 19).Value = 0.1443318412029361
$ws.Cells.Item(10, 20).Value = 0.1646881797369624

# Row 11: FAPs -> MuSCs
$ws.Cells.Item(11, 1).Value = "FAPs"
$ws.Cells.Item(11, 2).Value = "Wnt5a"
$ws.Cells.Item(11, 3).Value = "Fzd7"
$ws.Cells.Item(11, 4).Value = "MuSCs"
$ws.Cells.Item(11, 5).Value = 3
$ws.Cells.Item(11, 6).Value = 1
$ws.Cells.Item(11, 7).Value = 4.619088000000001
$ws.Cells.Item(11, 8).Value = 13.857264
$ws.Cells.Item(11, 9).Value = 0.9735252132746771
$ws.Cells.Item(11, 10).Value = 0.9821929966444328
$ws.Cells.Item(11, 11).Value = 2
$ws.Cells.Item(11, 12).Value = 1
$ws.Cells.Item(11, 13).Value = 6.825836
$ws.Cells.Item(11, 14).Value = 13.651672
$ws.Cells.Item(11, 15).Value = 0.2716022158275637
$ws.Cells.Item(11, 16).Value = 0.2047824327132465
$ws.Cells.Item(11, 17).Value = 31.529137157568
$ws.Cells.Item(11, 18).Value = 189.174822945408
$ws.Cells.Item(11, 19).Value = 0.2644116050894039
$ws.Cells.Item(11, 20).Value = 0.2011358712467606

# Row 12: FAPs -> Neutrophils
$ws.Cells.Item(12, 1).Value = "FAPs"
$ws.Cells.Item(12, 2).Value = "Wnt5a"
$ws.Cells.Item(12, 3).Value = "Fzd7"
$ws.Cells.Item(12, 4).Value = "Neutrophils"
$ws.Cells.Item(12, 5).Value = 3
$ws.Cells.Item(12, 6).Value = 1
$ws.Cells.Item(12, 7).Value = 4.619088000000001
$ws.Cells.Item(12, 8).Value = 13.857264
$ws.Cells.Item(12, 9).Value = 0.9735252132746771
$ws.Cells.Item(12, 10).Value = 0.9821929966444328
$ws.Cells.Item(12, 11).Value = 3
$ws.Cells.Item(12, 12).Value = 1
$ws.Cells.Item(12, 13).Value = 4.247626666666666
$ws.Cells.Item(12, 14).Value = 12.74288
$ws.Cells.Item(12, 15).Value = 0.1690144349607748
$ws.Cells.Item(12, 16).Value = 0.1911500632430207
$ws.Cells.Item(12, 17).Value = 19.62016136448
$ws.Cells.Item(12, 18).Value = 176.58145228032
$ws.Cells.Item(12, 19).Value = 0.1645398138416873
$ws.Cells.Item(12, 20).Value = 0.1877462534254354

# Row 13: FAPs -> Resolving-Mac
$ws.Cells.Item(13, 1).Value = "FAPs"
$ws.Cells.Item(13, 2).Value = "Wnt5a"
$ws.Cells.Item(13, 3).Value = "Fzd7"
$ws.Cells.Item(13, 4).Value = "Resolving-Mac"
$ws.Cells.Item(13, 5).Value = 3
$ws.Cells.Item(13, 6).Value = 1
$ws.Cells.Item(13, 7).Value = 4.619088000000001
$ws.Cells.Item(13, 8).Value = 13.857264
$ws.Cells.Item(13, 9).Value = 0.9735252132746771
$ws.Cells.Item(13, 10).Value = 0.9821929966444328
$ws.Cells.Item(13, 11).Value = 3
$ws.Cells.Item(13, 12).Value = 1
$ws.Cells.Item(13, 13).Value = 5.358931000000001
$ws.Cells.Item(13, 14).Value = 16.076793
$ws.Cells.Item(13, 15).Value = 0.213233592788784
$ws.Cells.Item(13, 16).Value = 0.2411605538696867
$ws.Cells.Item(13, 17).Value = 24.75337387492801
$ws.Cells.Item(13, 18).Value = 222.780364874352
$ws.Cells.Item(13, 19).Value = 0.2075882788970266
$ws.Cells.Item(13, 20).Value = 0.2368662070776987
